# TIMETABLE_LAST_WEEKS.xlsx update
# - Fill in previously-blank / differently-labelled "Mon/Tue/Wed/Thur" cells in the
#   "Fri (SWEN423 TEST)" table (T4:AA52) with "ENGR489 WRITING" for rows 16-25,
#   matching the style already used by neighbouring "ENGR489 WRITING" cells.
# - Fill in the "Fri" column (U47:U52) of that same table with "SWEN423 TEST",
#   matching the style already used higher up in the same column (U16).
# - Update the sheet view: zoom to 75%, and select T5:AA52 (the visible portion
#   of the third schedule table) as the last interacted-with range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant (used so the destination cell ends up with the exact
# same cell style / fill / font as the "source" neighbour cell).
$xlPasteFormats = -4122

function Copy-CellStyle-And-SetValue {
    param(
        [string]$TargetAddress,
        [string]$StyleSourceAddress,
        [string]$Value
    )

    $src = $ws.Range($StyleSourceAddress)
    $dst = $ws.Range($TargetAddress)

    # Copy just the formatting (fill colour / font / borders / etc.) from the
    # cell that already carries the "ENGR489 WRITING" / "SWEN423 TEST" style.
    $src.Copy()
    $dst.PasteSpecial($xlPasteFormats)

    # Now set the actual text value.
    $dst.Value = $Value
}

$excel.CutCopyMode = $false

# ---- "ENGR489 WRITING" block (Schedule5 table, columns X/Y/Z/AA) ----------

Copy-CellStyle-And-SetValue "X16"  "Y16"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X17"  "Y17"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X18"  "Y18"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X19"  "Y19"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X20"  "Y20"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "AA20" "Z20"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X21"  "Z21"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "Y21"  "Z21"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "AA21" "Z21"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X22"  "Z22"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "Y22"  "Z22"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "AA22" "Z22"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X23"  "V23"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "Y23"  "V23"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "Z23"  "V23"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "AA23" "V23"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X24"  "V24"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "Y24"  "V24"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "Z24"  "V24"  "ENGR489 WRITING"
Copy-CellStyle-And-SetValue "AA24" "V24"  "ENGR489 WRITING"

Copy-CellStyle-And-SetValue "X25"  "Y25"  "ENGR489 WRITING"

# ---- "SWEN423 TEST" block (Schedule5 table, "Fri" column U) ---------------

Copy-CellStyle-And-SetValue "U47" "U16" "SWEN423 TEST"
Copy-CellStyle-And-SetValue "U48" "U16" "SWEN423 TEST"
Copy-CellStyle-And-SetValue "U49" "U16" "SWEN423 TEST"
Copy-CellStyle-And-SetValue "U50" "U16" "SWEN423 TEST"
Copy-CellStyle-And-SetValue "U51" "U16" "SWEN423 TEST"
Copy-CellStyle-And-SetValue "U52" "U16" "SWEN423 TEST"

$excel.CutCopyMode = $false

# ---- Sheet view / selection -------------------------------------------------

$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 12
$win.Zoom = 75

$ws.Range("T5:AA52").Select()
